$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the per-row Price/Volume(1h) figures (and, for rows 45-46, the Coin
# name + Link too, since the two coins swapped rank) to match the latest
# coinranking.com snapshot.
#
# Every Price cell in the source file is stored as literal text, even when it
# happens to look like a plain number (e.g. "1.013"). Assigning such a string
# straight to Range.Value would get auto-coerced to a numeric value by COM, so
# those particular updates are written with a leading apostrophe (Excel's
# quote-prefix) to force text storage, and the default "Normal" style is then
# reapplied so the quote-prefix formatting left behind by that trick does not
# linger on the cell.

$ws.Range('D2').Value = '27.699.36'
$ws.Range('E2').Value = '  +2.22%  '
$ws.Range('D3').Value = '1.869.44'
$ws.Range('E3').Value = '  +1.04%  '
$ws.Range('D4').Value = "'1.013"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.27%  '
$ws.Range('D5').Value = "'313.28"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.21%  '
$ws.Range('E6').Value = '  -0.17%  '
$ws.Range('D7').Value = "'0.4817"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.04%  '
$ws.Range('D8').Value = "'0.3814"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.45%  '
$ws.Range('D9').Value = "'0.07372"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.76%  '
$ws.Range('D10').Value = "'0.9373"
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').Value = "'21.03"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +5.67%  '
$ws.Range('D12').Value = "'0.07801"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.08%  '
$ws.Range('D13').Value = '1.898.11'
$ws.Range('E13').Value = '  +3.29%  '
$ws.Range('D14').Value = "'5.484"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.61%  '
$ws.Range('D15').Value = "'6.587"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.58%  '
$ws.Range('D16').Value = "'90.57"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.45%  '
$ws.Range('E17').Value = '  -0.30%  '
$ws.Range('D18').Value = "'0.000008856"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.95%  '
$ws.Range('E19').Value = '  -0.18%  '
$ws.Range('D20').Value = '27.827.08'
$ws.Range('E20').Value = '  +2.67%  '
$ws.Range('D21').Value = "'14.75"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.13%  '
$ws.Range('D22').Value = "'5.122"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.26%  '
$ws.Range('D23').Value = '2.129.64'
$ws.Range('E23').Value = '  +2.75%  '
$ws.Range('D24').Value = "'10.79"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.33%  '
$ws.Range('D25').Value = "'1.944"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.02%  '
$ws.Range('D26').Value = "'156.60"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.40%  '
$ws.Range('E27').Value = '  +0.87%  '
$ws.Range('D28').Value = "'2.047"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.06%  '
$ws.Range('D29').Value = "'115.59"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.73%  '
$ws.Range('D30').Value = "'4.966"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.96%  '
$ws.Range('D31').Value = "'0.08915"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.53%  '
$ws.Range('E32').Value = '  +0.65%  '
$ws.Range('D33').Value = "'1.216"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.73%  '
$ws.Range('D34').Value = "'0.7612"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.16%  '
$ws.Range('D35').Value = "'4.628"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.21%  '
$ws.Range('D36').Value = "'2.729"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.50%  '
$ws.Range('D37').Value = "'1.133"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.65%  '
$ws.Range('D38').Value = "'0.02050"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.65%  '
$ws.Range('D39').Value = "'0.5676"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +7.54%  '
$ws.Range('D40').Value = "'0.05383"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.13%  '
$ws.Range('D41').Value = "'2.987"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.45%  '
$ws.Range('D42').Value = "'7.061"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.07%  '
$ws.Range('D43').Value = "'8.553"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.98%  '
$ws.Range('D44').Value = "'0.1531"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.46%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').Value = "'0.4907"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.58%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = "'10.74"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.19%  '
$ws.Range('E47').Value = '  -0.17%  '
$ws.Range('D48').Value = "'104.81"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.69%  '
$ws.Range('D49').Value = "'1.676"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.65%  '
$ws.Range('D50').Value = "'67.69"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.70%  '
$ws.Range('D51').Value = "'0.06104"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.69%  '
